# The document ends in a single empty paragraph right before the final
# sectPr. Replace it with two new paragraphs of plain text:
#   "Npm i"
#   "Npm run start"
# Each run is split the way Word splits text around words its spell
# checker flags ("Npm" and "start"), wrapping them in
# proofErr spellStart/spellEnd markers, matching the rest of the document.

$d = $word.ActiveDocument

$target = $d.Paragraphs.Last

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = '<w:p ' + $wNs + '>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Npm</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> i</w:t></w:r>' + `
  '</w:p>'

$para2 = '<w:p ' + $wNs + '>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Npm</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> run </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>start</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'

# Inserting this OOXML into the final (collapsed) paragraph's range
# overwrites that paragraph in place with the two new ones, instead of
# just being pushed in ahead of it.
[void]($target.Range.InsertXML($para1 + $para2))
